$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 8: financial period labels (shift left by one year, add new period) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Header row 9: publish dates (shift left by one, add new publish date) ---
$ws.Range("D9").Value = "1399-04-21 (13)"
$ws.Range("E9").Value = "1400-02-30 (7)"
$ws.Range("F9").Value = "1401-04-18 (7)"
$ws.Range("G9").Value = "1402-02-23 (8)"
$ws.Range("H9").Value = "1402-02-23"

# --- Data rows: shift each value one period to the left (D<-E<-F<-G<-H), new H is the new period ---
# Row 12
$ws.Range("D12").Value = 137572
$ws.Range("E12").Value = 503275
$ws.Range("F12").Value = 366420
$ws.Range("G12").Value = 1868143
$ws.Range("H12").Value = 5682077

# Row 13
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 318609
$ws.Range("G13").Value = 261339
$ws.Range("H13").Value = 346733

# Row 14
$ws.Range("D14").Value = 2757411
$ws.Range("E14").Value = 3751505
$ws.Range("F14").Value = 8320019
$ws.Range("G14").Value = 11169725
$ws.Range("H14").Value = 17358792

# Row 15
$ws.Range("D15").Value = 1135086
$ws.Range("E15").Value = 2162094
$ws.Range("F15").Value = 5060743
$ws.Range("G15").Value = 4608522
$ws.Range("H15").Value = 6407946

# Row 16
$ws.Range("D16").Value = 209726
$ws.Range("E16").Value = 1174656
$ws.Range("F16").Value = 1067892
$ws.Range("G16").Value = 1409203
$ws.Range("H16").Value = 2550680

# Row 17
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0

# Row 18
$ws.Range("D18").Value = 4239795
$ws.Range("E18").Value = 7591530
$ws.Range("F18").Value = 15133683
$ws.Range("G18").Value = 19316932
$ws.Range("H18").Value = 32346228

# Row 19
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0

# Row 20
$ws.Range("D20").Value = 152825
$ws.Range("E20").Value = 241772
$ws.Range("F20").Value = 456869
$ws.Range("G20").Value = 1146909
$ws.Range("H20").Value = 3997872

# Row 21
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 866660

# Row 22
$ws.Range("D22").Value = 778259
$ws.Range("E22").Value = 1168078
$ws.Range("F22").Value = 2378936
$ws.Range("G22").Value = 3300378
$ws.Range("H22").Value = 3581272

# Row 23
$ws.Range("D23").Value = 12032
$ws.Range("E23").Value = 16182
$ws.Range("F23").Value = 26704
$ws.Range("G23").Value = 27944
$ws.Range("H23").Value = 27462

# Row 24
$ws.Range("D24").Value = "-"
$ws.Range("E24").Value = "-"
$ws.Range("F24").Value = "-"
$ws.Range("G24").Value = "-"
$ws.Range("H24").Value = "-"

# Row 25
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0

# Row 26
$ws.Range("D26").Value = 943116
$ws.Range("E26").Value = 1426032
$ws.Range("F26").Value = 2862509
$ws.Range("G26").Value = 4475231
$ws.Range("H26").Value = 8473266

# Row 27
$ws.Range("D27").Value = 5182911
$ws.Range("E27").Value = 9017562
$ws.Range("F27").Value = 17996192
$ws.Range("G27").Value = 23792163
$ws.Range("H27").Value = 40819494

# Row 29
$ws.Range("D29").Value = 737023
$ws.Range("E29").Value = 1306229
$ws.Range("F29").Value = 4479869
$ws.Range("G29").Value = 5072864
$ws.Range("H29").Value = 10977719

# Row 30
$ws.Range("D30").Value = "-"
$ws.Range("E30").Value = "-"
$ws.Range("F30").Value = "-"
$ws.Range("G30").Value = "-"
$ws.Range("H30").Value = "-"

# Row 31
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0

# Row 32
$ws.Range("D32").Value = 133541
$ws.Range("E32").Value = 161036
$ws.Range("F32").Value = 155529
$ws.Range("G32").Value = 102254
$ws.Range("H32").Value = 52040

# Row 33
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = 991117
$ws.Range("G33").Value = 311296
$ws.Range("H33").Value = 820132

# Row 34
$ws.Range("D34").Value = 1717231
$ws.Range("E34").Value = 3820628
$ws.Range("F34").Value = 6577241
$ws.Range("G34").Value = 8302390
$ws.Range("H34").Value = 9180918

# Row 35
$ws.Range("D35").Value = 65791
$ws.Range("E35").Value = 42935
$ws.Range("F35").Value = 144435
$ws.Range("G35").Value = 222012
$ws.Range("H35").Value = 421647

# Row 36
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 0

# Row 37
$ws.Range("D37").Value = 2653586
$ws.Range("E37").Value = 5330828
$ws.Range("F37").Value = 12348191
$ws.Range("G37").Value = 14010816
$ws.Range("H37").Value = 21461706

# Row 38
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 0

# Row 39
$ws.Range("D39").Value = "-"
$ws.Range("E39").Value = "-"
$ws.Range("F39").Value = "-"
$ws.Range("G39").Value = "-"
$ws.Range("H39").Value = "-"

# Row 40
$ws.Range("D40").Value = 442000
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 2361111
$ws.Range("H40").Value = 7227166

# Row 41
$ws.Range("D41").Value = 38636
$ws.Range("E41").Value = 58178
$ws.Range("F41").Value = 87366
$ws.Range("G41").Value = 139607
$ws.Range("H41").Value = 243677

# Row 42
$ws.Range("D42").Value = 480636
$ws.Range("E42").Value = 58178
$ws.Range("F42").Value = 87366
$ws.Range("G42").Value = 2500718
$ws.Range("H42").Value = 7470843

# Row 43
$ws.Range("D43").Value = 3134222
$ws.Range("E43").Value = 5389006
$ws.Range("F43").Value = 12435557
$ws.Range("G43").Value = 16511534
$ws.Range("H43").Value = 28932549

# Row 45
$ws.Range("D45").Value = 600000
$ws.Range("E45").Value = 1000000
$ws.Range("F45").Value = 1000000
$ws.Range("G45").Value = 2000000
$ws.Range("H45").Value = 4000000

# Row 46
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0

# Row 47
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 0

# Row 48
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 0

# Row 49
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0

# Row 50
$ws.Range("D50").Value = 60000
$ws.Range("E50").Value = 100000
$ws.Range("F50").Value = 100000
$ws.Range("G50").Value = 200000
$ws.Range("H50").Value = 400000

# Row 51
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0

# Row 52
$ws.Range("D52").Value = "-"
$ws.Range("E52").Value = "-"
$ws.Range("F52").Value = "-"
$ws.Range("G52").Value = "-"
$ws.Range("H52").Value = "-"

# Row 53
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0

# Row 54
$ws.Range("D54").Value = "-"
$ws.Range("E54").Value = "-"
$ws.Range("F54").Value = "-"
$ws.Range("G54").Value = "-"
$ws.Range("H54").Value = "-"

# Row 55
$ws.Range("D55").Value = 0
$ws.Range("E55").Value = 0
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 0

# Row 56
$ws.Range("D56").Value = 1388689
$ws.Range("E56").Value = 2528556
$ws.Range("F56").Value = 4460635
$ws.Range("G56").Value = 5080629
$ws.Range("H56").Value = 7486945

# Row 57
$ws.Range("D57").Value = 2048689
$ws.Range("E57").Value = 3628556
$ws.Range("F57").Value = 5560635
$ws.Range("G57").Value = 7280629
$ws.Range("H57").Value = 11886945

# Row 58
$ws.Range("D58").Value = 5182911
$ws.Range("E58").Value = 9017562
$ws.Range("F58").Value = 17996192
$ws.Range("G58").Value = 23792163
$ws.Range("H58").Value = 40819494
